$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents (header text + values) of columns B and C for rows 1-5
for ($r = 1; $r -le 5; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $cVal = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 2).Value2 = $cVal
    $ws.Cells.Item($r, 3).Value2 = $bVal
}
